$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: columns A (Sending cluster), B (Ligand), C (Receptor), D (Target cluster), E..T metrics
$data = @(
    @("ECs", "Timp2", "Itga3", "ECs", 3, 1, 15.89290866666667, 47.678726, 0.04140398574586633, 0.04140398574586632, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 94.896446398306, 854.068017584754, 0.02549341749811079, 0.02549341749811078),
    @("ECs", "Timp2", "Itga3", "FAPs", 3, 1, 15.89290866666667, 47.678726, 0.04140398574586633, 0.04140398574586632, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 16.874263991016, 151.868375919144, 0.004533179831526202, 0.004533179831526201),
    @("ECs", "Timp2", "Itga3", "M2", 3, 1, 15.89290866666667, 47.678726, 0.04140398574586633, 0.04140398574586632, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 0.4691692591124444, 4.222523332012, 0.0001260397860382512, 0.0001260397860382512),
    @("ECs", "Timp2", "Itga3", "sCs", 3, 1, 15.89290866666667, 47.678726, 0.04140398574586633, 0.04140398574586632, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 41.88190940946644, 376.937184685198, 0.01125134863019109, 0.01125134863019109),
    @("FAPs", "Timp2", "Itga3", "ECs", 3, 1, 268.8003336666666, 806.401001, 0.700274909838328, 0.7002749098383279, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 1605.004910721331, 14445.04419649198, 0.4311758957105409, 0.4311758957105408),
    @("FAPs", "Timp2", "Itga3", "FAPs", 3, 1, 268.8003336666666, 806.401001, 0.700274909838328, 0.7002749098383279, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 285.3982166699159, 2568.583950029244, 0.07667068859716891, 0.0766706885971689),
    @("FAPs", "Timp2", "Itga3", "M2", 3, 1, 268.8003336666666, 806.401001, 0.700274909838328, 0.7002749098383279, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 7.935165050062444, 71.41648545056199, 0.002131739208532368, 0.002131739208532367),
    @("FAPs", "Timp2", "Itga3", "sCs", 3, 1, 268.8003336666666, 806.401001, 0.700274909838328, 0.7002749098383279, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 708.3581400976414, 6375.223260878773, 0.1902965863220858, 0.1902965863220857),
    @("M2", "Timp2", "Itga3", "ECs", 3, 1, 41.83092566666667, 125.492777, 0.1089773487260793, 0.1089773487260792, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 249.772164339187, 2247.949479052683, 0.06709994216410722, 0.0670999421641072),
    @("M2", "Timp2", "Itga3", "FAPs", 3, 1, 41.83092566666667, 125.492777, 0.1089773487260793, 0.1089773487260792, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 44.41390166473199, 399.725114982588, 0.01193155466651133, 0.01193155466651133),
    @("M2", "Timp2", "Itga3", "M2", 3, 1, 41.83092566666667, 125.492777, 0.1089773487260793, 0.1089773487260792, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 1.234876812963778, 11.113891316674, 0.000331742982445168, 0.0003317429824451679),
    @("M2", "Timp2", "Itga3", "sCs", 3, 1, 41.83092566666667, 125.492777, 0.1089773487260793, 0.1089773487260792, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 110.2352675668468, 992.1174081016209, 0.02961410891301556, 0.02961410891301554),
    @("sCs", "Timp2", "Itga3", "ECs", 3, 1, 57.32556, 171.97668, 0.1493437556897266, 0.1493437556897266, 2, 0.6666666666666666, 5.970993, 17.912979, 0.6157237531330177, 0.6157237531330177, 342.29051748108, 3080.61465732972, 0.09195449776025894, 0.09195449776025892),
    @("sCs", "Timp2", "Itga3", "FAPs", 3, 1, 57.32556, 171.97668, 0.1493437556897266, 0.1493437556897266, 3, 1, 1.061748, 3.185244, 0.10948655666511, 0.10948655666511, 60.86529867887999, 547.7876881099199, 0.0163511335699036, 0.0163511335699036),
    @("sCs", "Timp2", "Itga3", "M2", 3, 1, 57.32556, 171.97668, 0.1493437556897266, 0.1493437556897266, 1, 0.3333333333333333, 0.02952066666666667, 0.088562, 0.003044146203987976, 0.003044146203987975, 1.69228874824, 15.23059873416, 0.0004546242269721889, 0.0004546242269721887),
    @("sCs", "Timp2", "Itga3", "sCs", 3, 1, 57.32556, 171.97668, 0.1493437556897266, 0.1493437556897266, 3, 1, 2.635257666666666, 7.905773, 0.2717455439978843, 0.2717455439978843, 151.06762148596, 1359.60859337364, 0.04058350013259188, 0.04058350013259187),
)

$stringCols = @(0, 1, 2, 3)  # A, B, C, D are text columns (0-based index in row array)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($r, $col)
        if ($stringCols -contains $j) {
            $cell.Value = $row[$j]
        } else {
            $cell.Value2 = $row[$j]
        }
    }
}
